$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency prices in column D (stored as text values).
# NumberFormat "@" forces Excel to keep the numeric-looking string as text
# instead of auto-converting it to a number; ClearFormats() afterwards removes
# the now-unneeded explicit style so the cell stays styled like the original
# (default/general style), matching the source data which has no style index.
$priceUpdates = @{
    "D2" = "244.61"
    "D3" = "23.92"
    "D4" = "5.190"
    "D6" = "6.476"
    "D7" = "3.171"
    "D8" = "0.8137"
    "D9" = "0.8686"
    "D11" = "0.1372"
    "D12" = "0.06941"
    "D13" = "0.03177"
    "D14" = "0.02917"
    "D15" = "0.09323"
    "D16" = "3.827"
    "D17" = "0.001519"
    "D18" = "0.04724"
    "D19" = "0.006173"
    "D20" = "0.001236"
    "D21" = "0.004108"
    "D22" = "0.00008506"
    "D23" = "3.544"
    "D24" = "2.162"
    "D40" = "0.03695"
    "D41" = "0.006425"
    "D42" = "0.1051"
    "D44" = "0.008117"
    "D45" = "0.00005471"
    "D47" = "0.4003"
    "D48" = "0.002557"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cellRef]
    $range.ClearFormats()
}
